$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '(Intercept)'
$ws.Range('E2').Value = -1.23854976546635
$ws.Range('F2').Value = 0.0696205716433035
$ws.Range('G2').Value = -17.7899970688545
$ws.Range('H2').Value = 0.0000000000000000000000000000000000000000000000000000000000000000000000844843248705405

$ws.Range('D3').Value = 'habitat_typeExposed/Low SAV'
$ws.Range('E3').Value = 0.5907616369183
$ws.Range('F3').Value = 0.200489652741129
$ws.Range('G3').Value = 2.94659414509081
$ws.Range('H3').Value = 0.00321294597387518

$ws.Range('D4').Value = 'habitat_typeMod/Dense SAV'
$ws.Range('E4').Value = 0.925588437744798
$ws.Range('F4').Value = 0.383018311806287
$ws.Range('G4').Value = 2.41656445452905
$ws.Range('H4').Value = 0.0156677520554825

$ws.Range('D5').Value = 'habitat_typeShallow/Dense SAV'
$ws.Range('E5').Value = -0.14850617273924
$ws.Range('F5').Value = 0.0926687162520895
$ws.Range('G5').Value = -1.60254915299846
$ws.Range('H5').Value = 0.10903422734063

$ws.Range('D6').Value = 'habitat_typeShallow/Low SAV'
$ws.Range('E6').Value = 0.233126748906077
$ws.Range('F6').Value = 0.0216718930447652
$ws.Range('G6').Value = 10.7571013028042
$ws.Range('H6').Value = 0.00000000000000000000000000548689824848472

$ws.Range('D7').Value = 'seasonWinter'
$ws.Range('E7').Value = -0.239431936597665
$ws.Range('F7').Value = 0.0366427178206379
$ws.Range('G7').Value = -6.5342297416272
$ws.Range('H7').Value = 0.0000000000639377356070636

$ws.Range('D8').Value = 'seasonSpring'
$ws.Range('E8').Value = 0.769476436454697
$ws.Range('F8').Value = 0.152049489535731
$ws.Range('G8').Value = 5.0606972690551
$ws.Range('H8').Value = 0.000000417726013009785

$ws.Range('D9').Value = 'seasonSummer'
$ws.Range('E9').Value = 0.941956414450631
$ws.Range('F9').Value = 0.0221266307325715
$ws.Range('G9').Value = 42.5711634923262
$ws.Range('H9').Value = 0

$ws.Range('D10').Value = 'habitat_typeExposed/Low SAV:seasonWinter'
$ws.Range('E10').ClearContents()
$ws.Range('F10').ClearContents()
$ws.Range('G10').ClearContents()
$ws.Range('H10').ClearContents()

$ws.Range('D11').Value = 'habitat_typeMod/Dense SAV:seasonWinter'
$ws.Range('E11').ClearContents()
$ws.Range('F11').ClearContents()
$ws.Range('G11').ClearContents()
$ws.Range('H11').ClearContents()

$ws.Range('D12').Value = 'habitat_typeShallow/Dense SAV:seasonWinter'
$ws.Range('E12').Value = -0.300353534072304
$ws.Range('F12').Value = 0.047135774592608
$ws.Range('G12').Value = -6.37209288843228
$ws.Range('H12').Value = 0.000000000186465887056438

$ws.Range('D13').Value = 'habitat_typeShallow/Low SAV:seasonWinter'
$ws.Range('E13').Value = -0.358715774094696
$ws.Range('F13').Value = 0.0427597988870297
$ws.Range('G13').Value = -8.38908936504621
$ws.Range('H13').Value = 0.0000000000000000489922192413372

$ws.Range('D14').Value = 'habitat_typeExposed/Low SAV:seasonSpring'
$ws.Range('E14').Value = -0.544217780925337
$ws.Range('F14').Value = 0.319889208062018
$ws.Range('G14').Value = -1.70126958712476
$ws.Range('H14').Value = 0.0888923761087084

$ws.Range('D15').Value = 'habitat_typeMod/Dense SAV:seasonSpring'
$ws.Range('E15').Value = -0.784634891508651
$ws.Range('F15').Value = 0.443432840543253
$ws.Range('G15').Value = -1.76945598018268
$ws.Range('H15').Value = 0.0768178106557114

$ws.Range('D16').Value = 'habitat_typeShallow/Dense SAV:seasonSpring'
$ws.Range('E16').Value = 0.00536011573658114
$ws.Range('F16').Value = 0.157033280891721
$ws.Range('G16').Value = 0.0341336289106581
$ws.Range('H16').Value = 0.972770592108823

$ws.Range('D17').Value = 'habitat_typeShallow/Low SAV:seasonSpring'
$ws.Range('E17').Value = -0.309461459540699
$ws.Range('F17').Value = 0.157380828825096
$ws.Range('G17').Value = -1.96632246666216
$ws.Range('H17').Value = 0.0492613732191092

$ws.Range('D18').Value = 'habitat_typeExposed/Low SAV:seasonSummer'
$ws.Range('E18').Value = -0.592272478619254
$ws.Range('F18').Value = 0.232609096900932
$ws.Range('G18').Value = -2.54621374017673
$ws.Range('H18').Value = 0.0108898487026579

$ws.Range('D19').Value = 'habitat_typeMod/Dense SAV:seasonSummer'
$ws.Range('E19').Value = -0.948150275482377
$ws.Range('F19').Value = 0.395753442757862
$ws.Range('G19').Value = -2.39581055536766
$ws.Range('H19').Value = 0.0165836590905726

$ws.Range('D20').Value = 'habitat_typeShallow/Dense SAV:seasonSummer'
$ws.Range('E20').Value = -0.00105627700287568
$ws.Range('F20').Value = 0.0445431552978853
$ws.Range('G20').Value = -0.0237135648745977
$ws.Range('H20').Value = 0.981081085843704

$ws.Range('D21').Value = 'habitat_typeShallow/Low SAV:seasonSummer'
$ws.Range('E21').Value = -0.349663006339907
$ws.Range('F21').Value = 0.0313162617192958
$ws.Range('G21').Value = -11.1655410685388
$ws.Range('H21').Value = 0.0000000000000000000000000000601239954249567

$ws.Range('D22').Value = 'sd__(Intercept)'
$ws.Range('E22').Value = 0.189898454364348
$ws.Range('F22').ClearContents()
$ws.Range('G22').ClearContents()
$ws.Range('H22').ClearContents()
